$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '314.09'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2.11%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '40.91'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.30%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.155'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.48%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07595'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.95%'
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = 'FTXToken'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.680'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2.27%'
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9301'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.52%'
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.424'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.82%'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1200'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-3.83%'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1823'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.27%'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09059'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.78%'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04145'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-2.68%'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1053'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.22%'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001296'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '2.60%'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005832'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.34%'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'UpBots'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.007522'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.18%'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.331'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.51%'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.328'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.38%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.67%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.622'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-3.04%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '4.47%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04027'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.93%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001279'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.27%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.003974'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-7.32%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '2.21%'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-2.60%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05169'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-2.43%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007716'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.67%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1300'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.07%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007610'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '10.61%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '72.52%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008577'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '12.21%'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '10.79%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006594'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-2.05%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.17%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.2753'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '62.08%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '35.25%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002102'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.17%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002002'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.17%'
